# Generate Report for Handback
# This script updates the localization-status workbook to reflect that the
# de-de and zh-cn handback files have been generated: it updates the
# "Status" column, fills in "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" for the handback row on each language sheet,
# and widens a few columns so the new/longer content is readable.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Column width adjustments
# ---------------------------------------------------------------------

$wsOverview = $wb.Worksheets.Item("Overview")
# E (zh-cn) and F (de-de) columns widen from ~17.2 to ~30 characters
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$langSheets = @("zh-cn", "de-de")

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    # Column C (Status) widens from ~17.2 to ~30 characters
    $ws.Columns.Item(3).ColumnWidth = 29.166666666666668
    # Columns I (Latest Target File) and J (Latest Handback File) widen to 40
    $ws.Columns.Item(9).ColumnWidth = 39.166666666666664
    $ws.Columns.Item(10).ColumnWidth = 39.166666666666664
}

# ---------------------------------------------------------------------
# 2. Update the "Status" text for both language sheets
# ---------------------------------------------------------------------

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("C2").Value = "Handed back: in sync with en-US"
}

# ---------------------------------------------------------------------
# 3. zh-cn sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for the handback row
# ---------------------------------------------------------------------

$wsZh = $wb.Worksheets.Item("zh-cn")

# "Latest Target File" (I2) now references the same handed-back markdown
# file as the source file (A2), and becomes a hyperlink like A2.
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb387db6c93e3c0fd78d8b0cfb65a1901f2b5a01/e2e/fa4b0bfc-2c91-4de3-a3ea-50a815c75e2c.md", "", "", "fa4b0bfc-2c91-4de3-a3ea-50a815c75e2c.md")

# "Latest Handback File" (J2) is the latest generated xliff file (same as
# the "Latest Handoff File" in G2 for this sheet).
$wsZh.Range("J2").Value = "fa4b0bfc-2c91-4de3-a3ea-50a815c75e2c.47f35ddbea6e8d32461959e49b11bcec959dbc7f.zh-cn.xlf"

# "Latest Handback DateTime" (K2)
$wsZh.Range("K2").Value = "2016-09-05 11:38:39"

# ---------------------------------------------------------------------
# 4. de-de sheet: fill in Latest Target File / Latest Handback File /
#    Latest Handback DateTime for the handback row
# ---------------------------------------------------------------------

$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/eb387db6c93e3c0fd78d8b0cfb65a1901f2b5a01/e2e/fa4b0bfc-2c91-4de3-a3ea-50a815c75e2c.md", "", "", "fa4b0bfc-2c91-4de3-a3ea-50a815c75e2c.md")

$wsDe.Range("J2").Value = "fa4b0bfc-2c91-4de3-a3ea-50a815c75e2c.47f35ddbea6e8d32461959e49b11bcec959dbc7f.de-de.xlf"

$wsDe.Range("K2").Value = "2016-09-05 11:38:58"
